# Applies the weekly Fruta/hortaliza price-sheet update for
# "Fruta, Feria Lagunitas de Puerto Montt - Damasco": rows 28-54 are
# re-shuffled with refreshed values and two new records are appended
# as rows 55-56 (dimension grows from A1:T54 to A1:T56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("D28").Value = 44907
$ws.Range("K28").Value = 'Castle Brite'
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 400
$ws.Range("O28").Value = 22000
$ws.Range("P28").Value = 21500
$ws.Range("Q28").Value = '$/caja 16 kilos'
$ws.Range("R28").Value = 'Región de O''Higgins'
$ws.Range("S28").Value = 1344
$ws.Range("T28").Value = 16
# Row 29
$ws.Range("D29").Value = 44907
$ws.Range("K29").Value = 'Castle Brite'
$ws.Range("L29").Value = 'Segunda'
$ws.Range("Q29").Value = '$/caja 16 kilos'
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 1125
$ws.Range("T29").Value = 16
# Row 30
$ws.Range("L30").Value = 'Especial'
$ws.Range("N30").Value = 21000
$ws.Range("O30").Value = 21000
$ws.Range("P30").Value = 21000
$ws.Range("S30").Value = 1167
# Row 31
$ws.Range("D31").Value = 44575
$ws.Range("K31").Value = 'Modesto'
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 200
$ws.Range("N31").Value = 18000
$ws.Range("O31").Value = 18000
$ws.Range("P31").Value = 18000
$ws.Range("Q31").Value = '$/caja 18 kilos'
$ws.Range("S31").Value = 1000
$ws.Range("T31").Value = 18
# Row 32
$ws.Range("D32").Value = 44575
$ws.Range("K32").Value = 'Modesto'
$ws.Range("L32").Value = 'Segunda'
$ws.Range("N32").Value = 16000
$ws.Range("O32").Value = 16000
$ws.Range("P32").Value = 16000
$ws.Range("S32").Value = 889
# Row 33
$ws.Range("D33").Value = 44176
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 500
$ws.Range("N33").Value = 15000
$ws.Range("O33").Value = 16000
$ws.Range("P33").Value = 15500
$ws.Range("Q33").Value = '$/caja 15 kilos'
$ws.Range("S33").Value = 1033
$ws.Range("T33").Value = 15
# Row 34
$ws.Range("L34").Value = 'Especial'
$ws.Range("N34").Value = 20000
$ws.Range("O34").Value = 20000
$ws.Range("P34").Value = 20000
$ws.Range("S34").Value = 1111
# Row 35
$ws.Range("D35").Value = 44551
$ws.Range("M35").Value = 200
$ws.Range("N35").Value = 18000
$ws.Range("O35").Value = 18000
$ws.Range("P35").Value = 18000
$ws.Range("Q35").Value = '$/caja 18 kilos'
$ws.Range("S35").Value = 1000
$ws.Range("T35").Value = 18
# Row 36
$ws.Range("D36").Value = 44551
$ws.Range("L36").Value = 'Segunda'
$ws.Range("M36").Value = 200
$ws.Range("N36").Value = 16000
$ws.Range("O36").Value = 16000
$ws.Range("P36").Value = 16000
$ws.Range("S36").Value = 889
# Row 37
$ws.Range("D37").Value = 44894
$ws.Range("M37").Value = 400
$ws.Range("N37").Value = 23000
$ws.Range("O37").Value = 24000
$ws.Range("P37").Value = 23500
$ws.Range("Q37").Value = '$/caja 16 kilos'
$ws.Range("S37").Value = 1469
$ws.Range("T37").Value = 16
# Row 38
$ws.Range("L38").Value = 'Especial'
$ws.Range("N38").Value = 20000
$ws.Range("O38").Value = 20000
$ws.Range("P38").Value = 20000
$ws.Range("S38").Value = 1111
# Row 39
$ws.Range("D39").Value = 44547
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 350
$ws.Range("N39").Value = 18000
$ws.Range("O39").Value = 18000
$ws.Range("P39").Value = 18000
$ws.Range("S39").Value = 1000
# Row 40
$ws.Range("D40").Value = 44547
$ws.Range("L40").Value = 'Segunda'
$ws.Range("M40").Value = 350
$ws.Range("N40").Value = 16000
$ws.Range("O40").Value = 16000
$ws.Range("P40").Value = 16000
$ws.Range("S40").Value = 889
# Row 41
$ws.Range("L41").Value = 'Especial'
$ws.Range("N41").Value = 21000
$ws.Range("O41").Value = 21000
$ws.Range("P41").Value = 21000
$ws.Range("S41").Value = 1167
# Row 42
$ws.Range("D42").Value = 44568
$ws.Range("L42").Value = 'Primera'
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 18000
$ws.Range("O42").Value = 18000
$ws.Range("P42").Value = 18000
$ws.Range("Q42").Value = '$/caja 18 kilos'
$ws.Range("R42").Value = 'Región Metropolitana'
$ws.Range("S42").Value = 1000
$ws.Range("T42").Value = 18
# Row 43
$ws.Range("D43").Value = 44568
$ws.Range("L43").Value = 'Segunda'
$ws.Range("N43").Value = 16000
$ws.Range("O43").Value = 16000
$ws.Range("P43").Value = 16000
$ws.Range("Q43").Value = '$/caja 18 kilos'
$ws.Range("S43").Value = 889
$ws.Range("T43").Value = 18
# Row 44
$ws.Range("D44").Value = 44159
$ws.Range("L44").Value = 'Tercera'
$ws.Range("M44").Value = 400
$ws.Range("N44").Value = 15500
$ws.Range("O44").Value = 16000
$ws.Range("P44").Value = 15750
$ws.Range("Q44").Value = '$/caja 15 kilos'
$ws.Range("R44").Value = 'Región de O''Higgins'
$ws.Range("S44").Value = 1050
$ws.Range("T44").Value = 15
# Row 45
$ws.Range("D45").Value = 44900
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 200
$ws.Range("N45").Value = 23000
$ws.Range("O45").Value = 24000
$ws.Range("P45").Value = 23500
$ws.Range("Q45").Value = '$/caja 16 kilos'
$ws.Range("R45").Value = 'Región Metropolitana'
$ws.Range("S45").Value = 1469
$ws.Range("T45").Value = 16
# Row 46
$ws.Range("D46").Value = 44900
$ws.Range("M46").Value = 100
$ws.Range("N46").Value = 19000
$ws.Range("O46").Value = 19000
$ws.Range("P46").Value = 19000
$ws.Range("Q46").Value = '$/caja 16 kilos'
$ws.Range("R46").Value = 'Región Metropolitana'
$ws.Range("S46").Value = 1188
$ws.Range("T46").Value = 16
# Row 47
$ws.Range("D47").Value = 44530
$ws.Range("N47").Value = 20000
$ws.Range("O47").Value = 21000
$ws.Range("P47").Value = 20500
$ws.Range("R47").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S47").Value = 1139
# Row 48
$ws.Range("D48").Value = 44540
$ws.Range("L48").Value = 'Segunda'
$ws.Range("M48").Value = 600
$ws.Range("N48").Value = 16000
$ws.Range("O48").Value = 16000
$ws.Range("P48").Value = 16000
$ws.Range("R48").Value = 'Región del Maule'
$ws.Range("S48").Value = 889
# Row 49
$ws.Range("D49").Value = 44523
$ws.Range("L49").Value = 'Segunda'
$ws.Range("M49").Value = 500
$ws.Range("N49").Value = 28000
$ws.Range("O49").Value = 28500
$ws.Range("P49").Value = 28250
$ws.Range("R49").Value = 'Provincia de Limarí'
$ws.Range("S49").Value = 1569
# Row 50
$ws.Range("L50").Value = 'Especial'
$ws.Range("N50").Value = 20000
$ws.Range("O50").Value = 20000
$ws.Range("P50").Value = 20000
$ws.Range("S50").Value = 1111
# Row 51
$ws.Range("D51").Value = 44565
$ws.Range("M51").Value = 200
$ws.Range("N51").Value = 18000
$ws.Range("O51").Value = 18000
$ws.Range("P51").Value = 18000
$ws.Range("R51").Value = 'Región Metropolitana'
$ws.Range("S51").Value = 1000
# Row 52
$ws.Range("D52").Value = 44565
$ws.Range("M52").Value = 200
$ws.Range("N52").Value = 16000
$ws.Range("O52").Value = 16000
$ws.Range("P52").Value = 16000
$ws.Range("R52").Value = 'Región Metropolitana'
$ws.Range("S52").Value = 889
# Row 53
$ws.Range("D53").Value = 44537
$ws.Range("M53").Value = 500
$ws.Range("N53").Value = 20000
$ws.Range("O53").Value = 22000
$ws.Range("P53").Value = 21000
$ws.Range("R53").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S53").Value = 1167
# Row 54
$ws.Range("D54").Value = 44537
$ws.Range("M54").Value = 250
$ws.Range("N54").Value = 17000
$ws.Range("O54").Value = 17000
$ws.Range("P54").Value = 17000
$ws.Range("R54").Value = 'Región del Maule'
$ws.Range("S54").Value = 944
# Row 55 (new)
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C55").Value = 'Los Lagos'
$ws.Range("D55").Value = 44544
$ws.Range("E55").Value = 10
$ws.Range("F55").Value = 'Fruta'
$ws.Range("G55").Value = 100103
$ws.Range("H55").Value = 'Frutos de hueso (carozo)'
$ws.Range("I55").Value = 100103003
$ws.Range("J55").Value = 'Damasco'
$ws.Range("K55").Value = 'Castle Brite'
$ws.Range("L55").Value = 'Primera'
$ws.Range("M55").Value = 600
$ws.Range("N55").Value = 18000
$ws.Range("O55").Value = 20000
$ws.Range("P55").Value = 19000
$ws.Range("Q55").Value = '$/caja 18 kilos'
$ws.Range("R55").Value = 'Región Metropolitana'
$ws.Range("S55").Value = 1056
$ws.Range("T55").Value = 18
# Row 56 (new)
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C56").Value = 'Los Lagos'
$ws.Range("D56").Value = 44544
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 'Fruta'
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = 'Frutos de hueso (carozo)'
$ws.Range("I56").Value = 100103003
$ws.Range("J56").Value = 'Damasco'
$ws.Range("K56").Value = 'Castle Brite'
$ws.Range("L56").Value = 'Segunda'
$ws.Range("M56").Value = 300
$ws.Range("N56").Value = 16000
$ws.Range("O56").Value = 16000
$ws.Range("P56").Value = 16000
$ws.Range("Q56").Value = '$/caja 18 kilos'
$ws.Range("R56").Value = 'Región Metropolitana'
$ws.Range("S56").Value = 889
$ws.Range("T56").Value = 18

# New rows inherit the workbook's custom date style (cellXfs idx 2,
# numFmtId 165 'YYYY-MM-DD HH:MM:SS') on column D, matching every other
# data row; set explicitly since brand-new cells default to no style.
$ws.Range("D55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D56").NumberFormat = "YYYY-MM-DD HH:MM:SS"

